# Apply the new "word" (column C) sequence for rows 2-129 of the pairings sheet.
# Columns A (index), B (image), and D (category) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    "biegen", "lockern", "hören", "warnen", "schenken", "treiben",
    "schmecken", "zögern", "fällen", "sterben", "bauen", "siegen",
    "töten", "sichern", "seufzen", "stehlen", "sinken", "stammen",
    "lügen", "erben", "klingen", "schwingen", "wirken", "loben",
    "wüten", "fischen", "trennen", "trauen", "kümmern", "planen",
    "spinnen", "schreiten", "quälen", "werden", "kichern", "gelten",
    "klettern", "kehren", "tropfen", "rufen", "zielen", "bergen",
    "sprengen", "heben", "bitten", "greifen", "platzen", "münzen",
    "dienen", "wachsen", "ärgern", "heulen", "filmen", "arten",
    "enden", "liegen", "lesen", "fangen", "heilen", "grüßen",
    "weichen", "fahren", "betteln", "kosten", "dringen", "sperren",
    "mögen", "scheitern", "führen", "schlucken", "spielen", "ändern",
    "knarren", "jubeln", "schrecken", "jagen", "spüren", "flüchten",
    "drehen", "boxen", "streichen", "machen", "freuen", "suchen",
    "rasen", "sorgen", "feiern", "runden", "bluten", "folgen",
    "flehen", "scheinen", "räumen", "liefern", "äußern", "formen",
    "helfen", "decken", "tollen", "ehren", "schulden", "irren",
    "backen", "altern", "hauen", "schwächen", "reizen", "wehtun",
    "graben", "werfen", "fallen", "achten", "wenden", "gründen",
    "pflanzen", "mauern", "wundern", "schwören", "stecken", "geben",
    "pfeifen", "fließen", "zeigen", "zünden", "saufen", "malen",
    "bellen", "brauchen"
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $words[$i]
}
